{"js": "const replacements = {\n  \"27\\u00d746=\": \"73\\u00d728=\",\n  \"60\\u00d782=\": \"76\\u00d731=\",\n  \"76\\u00d747=\": \"67\\u00d783=\",\n  \"36\\u00d794=\": \"71\\u00d761=\",\n  \"58\\u00d773=\": \"28\\u00d739=\",\n  \"49\\u00d757=\": \"94\\u00d753=\",\n  \"98\\u00d782=\": \"37\\u00d781=\",\n  \"67\\u00d791=\": \"80\\u00d754=\",\n  \"50\\u00d724=\": \"20\\u00d765=\",\n  \"82\\u00d714=\": \"24\\u00d719=\",\n  \"85\\u00d717=\": \"35\\u00d750=\",\n  \"83\\u00d763=\": \"30\\u00d742=\",\n  \"78\\u00d713=\": \"19\\u00d735=\",\n  \"29\\u00d721=\": \"65\\u00d751=\",\n  \"96\\u00d788=\": \"68\\u00d797=\",\n  \"95\\u00d755=\": \"88\\u00d759=\",\n  \"22\\u00d739=\": \"94\\u00d786=\",\n  \"17\\u00d766=\": \"86\\u00d776=\",\n  \"92\\u00d779=\": \"67\\u00d765=\",\n  \"44\\u00d770=\": \"93\\u00d750=\",\n  \"25\\u00d732=\": \"37\\u00d772=\",\n  \"94\\u00d774=\": \"66\\u00d734=\",\n  \"68\\u00d799=\": \"81\\u00d745=\",\n  \"98\\u00d785=\": \"54\\u00d723=\",\n  \"35\\u00d714=\": \"62\\u00d715=\",\n};\n\nfor (const oldText of Object.keys(replacements)) {\n  const newText = replacements[oldText];\n  const searchResults = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  searchResults.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < searchResults.items.length; i++) {\n    searchResults.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# U+00D7 MULTIPLICATION SIGN, built from its code point to avoid any\n# source-encoding ambiguity for the literal glyph.\n$x = [char]0x00D7\n\n$replacements = @(\n  @(\"27${x}46=\", \"73${x}28=\"),\n  @(\"60${x}82=\", \"76${x}31=\"),\n  @(\"76${x}47=\", \"67${x}83=\"),\n  @(\"36${x}94=\", \"71${x}61=\"),\n  @(\"58${x}73=\", \"28${x}39=\"),\n  @(\"49${x}57=\", \"94${x}53=\"),\n  @(\"98${x}82=\", \"37${x}81=\"),\n  @(\"67${x}91=\", \"80${x}54=\"),\n  @(\"50${x}24=\", \"20${x}65=\"),\n  @(\"82${x}14=\", \"24${x}19=\"),\n  @(\"85${x}17=\", \"35${x}50=\"),\n  @(\"83${x}63=\", \"30${x}42=\"),\n  @(\"78${x}13=\", \"19${x}35=\"),\n  @(\"29${x}21=\", \"65${x}51=\"),\n  @(\"96${x}88=\", \"68${x}97=\"),\n  @(\"95${x}55=\", \"88${x}59=\"),\n  @(\"22${x}39=\", \"94${x}86=\"),\n  @(\"17${x}66=\", \"86${x}76=\"),\n  @(\"92${x}79=\", \"67${x}65=\"),\n  @(\"44${x}70=\", \"93${x}50=\"),\n  @(\"25${x}32=\", \"37${x}72=\"),\n  @(\"94${x}74=\", \"66${x}34=\"),\n  @(\"68${x}99=\", \"81${x}45=\"),\n  @(\"98${x}85=\", \"54${x}23=\"),\n  @(\"35${x}14=\", \"62${x}15=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
